$wb = $excel.ActiveWorkbook

# This script applies updated market-board price data (and recomputed
# profit columns) to specific Leve rows across all 8 job sheets, matching
# the scheduled-runner price refresh described in the commit diff.


# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 383800.34
$ws.Range("J17").Value = 383800.34
$ws.Range("L17").Value = 1151401.02
$ws.Range("N17").Value = -1151737.02

# Row 40
$ws.Range("H40").Value = 5810.5557
$ws.Range("J40").Value = 7500
$ws.Range("L40").Value = 7500
$ws.Range("N40").Value = -7850

# Row 86
$ws.Range("H86").Value = 2141.8823
$ws.Range("I86").Value = 1956.4546
$ws.Range("J86").Value = 2481.8333
$ws.Range("K86").Value = 1956.4546
$ws.Range("L86").Value = 2481.8333
$ws.Range("M86").Value = -833.4546
$ws.Range("N86").Value = -4727.8333

# Row 89
$ws.Range("H89").Value = 2141.8823
$ws.Range("I89").Value = 1956.4546
$ws.Range("J89").Value = 2481.8333
$ws.Range("K89").Value = 9782.273000000001
$ws.Range("L89").Value = 12409.1665
$ws.Range("M89").Value = -4166.273000000001
$ws.Range("N89").Value = -23641.1665

# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# Row 112
$ws.Range("H112").Value = 2161
$ws.Range("J112").Value = 2534.2856
$ws.Range("L112").Value = 7602.8568
$ws.Range("N112").Value = -9818.856800000001


# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2641.125
$ws.Range("I61").Value = 1297.3
$ws.Range("J61").Value = 4880.8335
$ws.Range("K61").Value = 1297.3
$ws.Range("L61").Value = 4880.8335
$ws.Range("M61").Value = -1085.3
$ws.Range("N61").Value = -5304.8335

# Row 122
$ws.Range("H122").Value = 2190.7188
$ws.Range("J122").Value = 3170.1428
$ws.Range("L122").Value = 9510.428400000001
$ws.Range("N122").Value = -14410.4284

# Row 132
$ws.Range("H132").Value = 2435.9456
$ws.Range("I132").Value = 1929.3617
$ws.Range("K132").Value = 5788.0851
$ws.Range("M132").Value = -3258.0851

# Row 136
$ws.Range("H136").Value = 2641.125
$ws.Range("I136").Value = 1297.3
$ws.Range("J136").Value = 4880.8335
$ws.Range("K136").Value = 3891.9
$ws.Range("L136").Value = 14642.5005
$ws.Range("M136").Value = -1341.9
$ws.Range("N136").Value = -19742.5005


# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 42
$ws.Range("H42").Value = 129997
$ws.Range("J42").Value = 129997
$ws.Range("L42").Value = 129997
$ws.Range("N42").Value = -130653

# Row 99
$ws.Range("H99").Value = 27507.625
$ws.Range("I99").Value = 31251.572
$ws.Range("K99").Value = 31251.572
$ws.Range("M99").Value = -29753.572

# Row 115
$ws.Range("H115").Value = 70000
$ws.Range("J115").Value = 70000
$ws.Range("L115").Value = 70000
$ws.Range("N115").Value = -73134

# Row 134
$ws.Range("H134").Value = 1561.1333
$ws.Range("I134").Value = 1561.1333
$ws.Range("K134").Value = 4683.3999
$ws.Range("M134").Value = -2148.3999


# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

# Row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

# Row 105
$ws.Range("H105").Value = 2470.818
$ws.Range("I105").Value = 2275.625
$ws.Range("J105").Value = 2991.3333
$ws.Range("K105").Value = 2275.625
$ws.Range("L105").Value = 2991.3333
$ws.Range("M105").Value = -528.625
$ws.Range("N105").Value = -6485.3333


# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Range("H56").Value = 10010.134
$ws.Range("I56").Value = 10010.134
$ws.Range("K56").Value = 10010.134
$ws.Range("M56").Value = -9480.134


# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3521.2727
$ws.Range("I80").Value = 2860.4443
$ws.Range("J80").Value = 6495
$ws.Range("K80").Value = 2860.4443
$ws.Range("L80").Value = 6495
$ws.Range("M80").Value = -1862.4443
$ws.Range("N80").Value = -8491

# Row 83
$ws.Range("H83").Value = 3521.2727
$ws.Range("I83").Value = 2860.4443
$ws.Range("J83").Value = 6495
$ws.Range("K83").Value = 14302.2215
$ws.Range("L83").Value = 32475
$ws.Range("M83").Value = -9310.2215
$ws.Range("N83").Value = -42459

# Row 87
$ws.Range("H87").Value = 23354
$ws.Range("J87").Value = 23354
$ws.Range("L87").Value = 23354
$ws.Range("N87").Value = -25850

# Row 90
$ws.Range("H90").Value = 23354
$ws.Range("J90").Value = 23354
$ws.Range("L90").Value = 70062
$ws.Range("N90").Value = -82542

# Row 102
$ws.Range("H102").Value = 1265.4517
$ws.Range("I102").Value = 1329.9642
$ws.Range("K102").Value = 1329.9642
$ws.Range("M102").Value = 292.0358000000001

# Row 122
$ws.Range("H122").Value = 3212.5806
$ws.Range("I122").Value = 3078.4644
$ws.Range("K122").Value = 9235.393199999999
$ws.Range("M122").Value = -6785.393199999999


# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 1999
$ws.Range("I68").Value = 1999
$ws.Range("K68").Value = 1999
$ws.Range("M68").Value = -1250

# Row 71
$ws.Range("H71").Value = 1999
$ws.Range("I71").Value = 1999
$ws.Range("K71").Value = 9995
$ws.Range("M71").Value = -6251

# Row 80
$ws.Range("H80").Value = 24790.666
$ws.Range("J80").Value = 25128
$ws.Range("L80").Value = 25128
$ws.Range("N80").Value = -27374

# Row 82
$ws.Range("H82").Value = 3785.25
$ws.Range("I82").Value = 3570.5
$ws.Range("K82").Value = 3570.5
$ws.Range("M82").Value = -3209.5

# Row 83
$ws.Range("H83").Value = 24790.666
$ws.Range("J83").Value = 25128
$ws.Range("L83").Value = 75384
$ws.Range("N83").Value = -86616

# Row 85
$ws.Range("H85").Value = 3785.25
$ws.Range("I85").Value = 3570.5
$ws.Range("K85").Value = 3570.5
$ws.Range("M85").Value = -2322.5

# Row 117
$ws.Range("H117").Value = 42000
$ws.Range("I117").Value = 42000
$ws.Range("K117").Value = 42000
$ws.Range("M117").Value = -37411


# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

# Row 62
$ws.Range("H62").Value = 14559.182
$ws.Range("I62").Value = 5018.875
$ws.Range("J62").Value = 40000
$ws.Range("K62").Value = 5018.875
$ws.Range("L62").Value = 40000
$ws.Range("M62").Value = -4394.875
$ws.Range("N62").Value = -41248

# Row 65
$ws.Range("H65").Value = 14559.182
$ws.Range("I65").Value = 5018.875
$ws.Range("J65").Value = 40000
$ws.Range("K65").Value = 25094.375
$ws.Range("L65").Value = 200000
$ws.Range("M65").Value = -21974.375
$ws.Range("N65").Value = -206240

# Row 68
$ws.Range("H68").Value = 25271
$ws.Range("J68").Value = 25271
$ws.Range("L68").Value = 25271
$ws.Range("N68").Value = -26893

# Row 69
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

# Row 71
$ws.Range("H71").Value = 25271
$ws.Range("J71").Value = 25271
$ws.Range("L71").Value = 75813
$ws.Range("N71").Value = -83925

# Row 72
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

# Row 81
$ws.Range("H81").Value = 2096.6667
$ws.Range("I81").Value = 2096.6667
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 4193.3334
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -3132.3334
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 2096.6667
$ws.Range("I84").Value = 2096.6667
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 20966.667
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -15662.667
$ws.Range("N84").ClearContents()

# Row 96
$ws.Range("H96").Value = 3910.4211
$ws.Range("I96").Value = 2300
$ws.Range("J96").Value = 3999.889
$ws.Range("K96").Value = 2300
$ws.Range("L96").Value = 3999.889
$ws.Range("M96").Value = -927
$ws.Range("N96").Value = -6745.889

# Row 136
$ws.Range("H136").Value = 2035.0526
$ws.Range("I136").Value = 1697.875
$ws.Range("K136").Value = 5093.625
$ws.Range("M136").Value = -2543.625
